$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.671.65'
$ws.Range("E2").Value = '  -2.51%  '
$ws.Range("D3").Value = '3.203.93'
$ws.Range("E3").Value = '  -3.26%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = "'593.90"
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("D6").Value = "'136.35"
$ws.Range("E6").Value = '  -5.40%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '3.202.09'
$ws.Range("E8").Value = '  -3.09%  '
$ws.Range("D9").Value = "'0.508"
$ws.Range("E9").Value = '  -3.07%  '
$ws.Range("E10").Value = '  -3.65%  '
$ws.Range("D11").Value = "'5.36"
$ws.Range("E11").Value = '  -2.43%  '
$ws.Range("E12").Value = '  -4.40%  '
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("E13").Value = '  -4.69%  '
$ws.Range("D14").Value = "'33.61"
$ws.Range("E14").Value = '  -4.46%  '
$ws.Range("D15").Value = '3.735.57'
$ws.Range("E15").Value = '  -3.12%  '
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '3.206.58'
$ws.Range("E17").Value = '  -3.12%  '
$ws.Range("D18").Value = '62.742.35'
$ws.Range("E18").Value = '  -2.50%  '
$ws.Range("D19").Value = "'6.71"
$ws.Range("E19").Value = '  -3.13%  '
$ws.Range("D20").Value = "'466.37"
$ws.Range("E20").Value = '  -4.25%  '
$ws.Range("D21").Value = "'13.91"
$ws.Range("E21").Value = '  -3.51%  '
$ws.Range("D22").Value = "'0.715"
$ws.Range("E22").Value = '  -4.26%  '
$ws.Range("D23").Value = "'7.70"
$ws.Range("E23").Value = '  -4.79%  '
$ws.Range("D24").Value = "'84.50"
$ws.Range("E24").Value = '  -0.78%  '
$ws.Range("D25").Value = "'13.40"
$ws.Range("E25").Value = '  -1.88%  '
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = '  -3.64%  '
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("D29").Value = "'7.90"
$ws.Range("E29").Value = '  -5.43%  '
$ws.Range("D30").Value = "'6.94"
$ws.Range("E30").Value = '  -4.57%  '
$ws.Range("D31").Value = "'2.09"
$ws.Range("E31").Value = '  -4.34%  '
$ws.Range("D32").Value = "'27.72"
$ws.Range("E32").Value = '  -2.74%  '
$ws.Range("E33").Value = '  -5.36%  '
$ws.Range("D34").Value = "'2.44"
$ws.Range("E34").Value = '  -6.00%  '
$ws.Range("E35").Value = '  -4.62%  '
$ws.Range("D36").Value = "'5.87"
$ws.Range("E36").Value = '  -2.92%  '
$ws.Range("D37").Value = "'51.76"
$ws.Range("E37").Value = '  -2.96%  '
$ws.Range("D38").Value = '0.0₃0702'
$ws.Range("E38").Value = '  -5.97%  '
$ws.Range("D39").Value = "'0.0392"
$ws.Range("E39").Value = '  -2.40%  '
$ws.Range("D40").Value = "'420.23"
$ws.Range("E40").Value = '  -2.69%  '
$ws.Range("D41").Value = '3.016.05'
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").Value = "'0.116"
$ws.Range("E42").Value = '  +4.72%  '
$ws.Range("D43").Value = "'8.12"
$ws.Range("E43").Value = '  -4.55%  '
$ws.Range("D44").Value = "'2.63"
$ws.Range("E44").Value = '  -6.83%  '
$ws.Range("E45").Value = '  -6.65%  '
$ws.Range("D46").Value = "'2.15"
$ws.Range("E46").Value = '  -4.96%  '
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("D48").Value = "'35.77"
$ws.Range("E48").Value = '  +2.38%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = "'25.68"
$ws.Range("E49").Value = '  -2.97%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = "'125.42"
$ws.Range("E50").Value = '  +1.46%  '
$ws.Range("E51").Value = '  -2.93%  '

# Reset style on cells forced to text via quote-prefix, so no stray
# cell-style index is introduced (matches original unstyled cells).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
